$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 values
$ws.Range("A1").Value = 0.76
$ws.Range("B1").Value = -0.18
$ws.Range("C1").Value = -0.17

# Clear previous D1 value (row now only spans A:C)
$ws.Range("D1").ClearContents()

# Add new row 2 values
$ws.Range("A2").Value = -0.8100000000000001
$ws.Range("B2").Value = -0.31
$ws.Range("C2").Value = -0.8100000000000001
